# B1--and-B2-PowerPoint.pptx edit
#
# 1) Slide 5's table switches to a different (built-in) table style.
# 2) The presentation's theme colour scheme switches from the "Red Violet"
#    (Integral) palette to the standard "Office" palette.

$p = $ppt.ActivePresentation

# --- 1) Update the table style on slide 5 -----------------------------
$slide5 = $p.Slides.Item(5)
for ($i = 1; $i -le $slide5.Shapes.Count; $i++) {
    $shp = $slide5.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{BA31744D-5BA4-4FF3-9328-FE44CB167932}")
    }
}

# --- 2) Swap the theme colour scheme to the Office palette -------------
# RRGGBB -> the BGR-ordered integer PowerPoint's ColorFormat.RGB expects.
function Convert-HexToBgr([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $b * 65536 + $g * 256 + $r
}

$officeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme

for ($i = 1; $i -le $officeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = Convert-HexToBgr $officeColors[$i - 1]
}
